# addressbook.xlsx: the header row (Name / Phone / E-Mail / Address) is
# replaced with the first real contact record written by the new
# "write to addressbook.xlsx" feature.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mark"

# Phone numbers can start with a leading zero ("010...") -- format the
# cell as Text first so Excel doesn't coerce the entry into a Number and
# silently drop the leading zero.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "01026134341"

$ws.Range("C1").Value = "hagi1126"
$ws.Range("D1").Value = "seoul"

[void]$ws.Range("A1").Select()
